# Weekly fruit/vegetable price update:
# Insert a new weekly record at row 133 (shifting the existing rows 133-225
# down to 134-226) for "Zapallo italiano" at Mercado Mayorista Lo Valledor
# de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 133:225 down to 134:226, leaving row 133 empty for the new record.
$ws.Rows("133:133").Insert()

# Populate the new row 133 with the new weekly observation.
$ws.Range("A133").Value = 6
$ws.Range("B133").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C133").Value = "Metropolitana"
$ws.Range("D133").Value = 44488
$ws.Range("E133").Value = 13
$ws.Range("F133").Value = 100112032
$ws.Range("G133").Value = "Zapallo italiano"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 420
$ws.Range("K133").Value = 13000
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = 14095
$ws.Range("N133").Value = "`$/caja 50 unidades"
$ws.Range("O133").Value = "Región de O'Higgins"
$ws.Range("P133").Value = 282
$ws.Range("Q133").Value = 50
$ws.Range("R133").Value = "Hortaliza"
